# ---------------------------------------------------------------------------
# Adds two new quiz sheets ("10_" and "11_") to the workbook and fills in
# content for three existing, previously-blank quiz sheets ("7_", "8_", "9_").
# Mirrors an "Add files via upload" commit that populated placeholder rows
# with question / answer / comment text for a set of True-False and
# multiple-select questions about internal (pipe) convection.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "9_" (Worksheets.Item(10)) is used as the template: copy it
# twice (while it is still blank) to create sheets "10_" and "11_"
# with identical column widths / styles / page setup, then rename.
# -----------------------------------------------------------------
$template = $wb.Worksheets.Item(10)

$template.Copy([System.Reflection.Missing]::Value, $template)
$ws10 = $wb.Worksheets.Item(11)
$ws10.Name = "10_"

$template.Copy([System.Reflection.Missing]::Value, $ws10)
$ws11 = $wb.Worksheets.Item(12)
$ws11.Name = "11_"

# -----------------------------------------------------------------
# Sheet "7_" -- True/False question about the temperature profile
# staying constant in the axial direction for fully developed flow.
# -----------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(8)
$ws7.Cells.Item(2,1).Value = "True or False: In fully developed internal flow, the the temperature profile within the flow stays the same in the axial direction (that is, along the length of the pipe)."
$ws7.Cells.Item(2,2).Value = "F"
$ws7.Cells.Item(2,3).Value = "The *shape* of the profile stays the same, but temperatures will actually continue to go up or down (depending on the surface condition)."
$ws7.Rows.Item(2).RowHeight = 105
$ws7.Range("A2:C2").Select()

# -----------------------------------------------------------------
# Sheet "8_" -- multi-select question about the implications of a
# constant positive surface flux on the temperature profile.
# -----------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(9)
$ws8.Cells.Item(1,1).Value = "The shape of the temperature profile stays the same in fully developed internal flow.   If the boundary condition is a constant positive flux (heat into the flow), what does this imply? Submit all that are true"
$ws8.Rows.Item(1).RowHeight = 135

$ws8.Cells.Item(2,1).Value = "dT/dr at the surface is constant"
$ws8.Cells.Item(2,2).Value = "Y"
$ws8.Cells.Item(2,3).Value = "If flux is constant, then the temperature gradient at the surface must remain the same."
$ws8.Rows.Item(2).RowHeight = 45

$ws8.Cells.Item(3,1).Value = "Surface temperature is constant"
$ws8.Cells.Item(3,2).Value = "N"
$ws8.Cells.Item(3,3).Value = "The mean temperature is going up at a constant rate (because energy is being added to the flow).  In order for that continue to happen, the surface temperature must go up."
$ws8.Rows.Item(3).RowHeight = 90

$ws8.Cells.Item(4,1).Value = "Mean temperature is moving up at a constant rate"
$ws8.Cells.Item(4,2).Value = "Y"
$ws8.Cells.Item(4,3).Value = "This is because a constant flux of energy is being added to the flow."
$ws8.Rows.Item(4).RowHeight = 45

$ws8.Cells.Item(5,1).Value = "The shape of the temperature profile remains exactly the same"
$ws8.Cells.Item(5,2).Value = "Y"
$ws8.Cells.Item(5,3).Value = "Not only is the mathematical shape of the profile the same, but the boundary conditions of the profile (dT/dr at the surfaces) is the same, so the profile shape stays exactly the same."
$ws8.Rows.Item(5).RowHeight = 90

$ws8.Range("C5").Select()

# -----------------------------------------------------------------
# Sheet "9_" -- True/False question about Prandtl number and the
# thermal boundary-layer entrance region.
# -----------------------------------------------------------------
$ws9 = $wb.Worksheets.Item(10)
$ws9.Cells.Item(2,1).Value = "True or False: A higher Prandtl number results in a shorter thermal boundary entrance region."
$ws9.Cells.Item(2,2).Value = "F"
$ws9.Cells.Item(2,3).Value = "A high Prandtl number implies that thermal diffusion is slower than momentum diffusion.  So the thermal boundary layer will grow more slowly, and the boundary layers from the surfaces won't meet in the center of the flow until farther into the pipe."
$ws9.Rows.Item(2).RowHeight = 120
$ws9.Range("A1:XFD1048576").Select()

# -----------------------------------------------------------------
# Sheet "10_" -- multiple choice: behavior of T(x) for equation 1
# (constant, positive flux boundary condition).
# -----------------------------------------------------------------
$ws10.Cells.Item(1,1).Value = "Look at equation 1 (dT/dx for a constant flux) on the last slide of the video.   If the flux is positive (into the pipe), what does the equation tell us about the temperature behavior over the length of the pipe? "
$ws10.Rows.Item(1).RowHeight = 120

$ws10.Cells.Item(2,1).Value = "The function T(x) will be constant"
$ws10.Cells.Item(2,2).Value = "N"
$ws10.Rows.Item(2).RowHeight = 30

$ws10.Cells.Item(3,1).Value = "The function T(x) will be linear"
$ws10.Cells.Item(3,2).Value = "Y"
$ws10.Cells.Item(3,3).Value = "The equation tells us that the slope of dT/dx is equal to a constant (everything on the right hand side is constant), so the temperature will increase linearly."
$ws10.Rows.Item(3).RowHeight = 90

$ws10.Cells.Item(4,1).Value = "The function T(x) will be a negative exponential function (increasing quickly and then more slowly)"
$ws10.Cells.Item(4,2).Value = "N"
$ws10.Rows.Item(4).RowHeight = 60

$ws10.Cells.Item(5,1).Value = "The function T(x) will be a positive exponential function (increasing slowly and then more quickly)"
$ws10.Cells.Item(5,2).Value = "N"
$ws10.Rows.Item(5).RowHeight = 60

$ws10.Range("A2:A5").Select()

# -----------------------------------------------------------------
# Sheet "11_" -- multiple choice: behavior of T(x) for equation 2
# (constant surface temperature boundary condition).
# -----------------------------------------------------------------
$ws11.Cells.Item(1,1).Value = "Look at equation 2 (dT/dx for a constant surface temperature) on the last slide of the video.   If the fluid temperature is initially cooler than the surface temperature , what does the equation tell us about the temperature behavior over the length of the pipe? "
$ws11.Rows.Item(1).RowHeight = 165

$ws11.Cells.Item(2,1).Value = "The function T(x) will be constant"
$ws11.Cells.Item(2,2).Value = "N"
$ws11.Rows.Item(2).RowHeight = 30

$ws11.Cells.Item(3,1).Value = "The function T(x) will be linear"
$ws11.Cells.Item(3,2).Value = "N"
$ws11.Rows.Item(3).RowHeight = 30

$ws11.Cells.Item(4,1).Value = "The function T(x) will be a negative exponential function (increasing quickly and then more slowly)"
$ws11.Cells.Item(4,2).Value = "Y"
$ws11.Cells.Item(4,3).Value = "The right hand side of the equation is not constant here: it will initially be large, but get smaller as the mean temperature approaches the surface temperature.  So the slope will be initially larger, but flatten out as the mean temp asymptotically approaches the surface temperature.  This is a negative exponential function."
$ws11.Rows.Item(4).RowHeight = 165

$ws11.Cells.Item(5,1).Value = "The function T(x) will be a positive exponential function (increasing slowly and then more quickly)"
$ws11.Cells.Item(5,2).Value = "N"
$ws11.Rows.Item(5).RowHeight = 60

$ws11.Range("D4").Select()

# "11_" is the sheet that ends up active / in front in the final workbook.
$ws11.Activate()
